$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting from the cell directly above (row 7) down into row 8
# so that C8/D8 pick up the same style (s="2") already used by C7/D7.
$ws.Range("C7:D7").Copy() | Out-Null
$ws.Range("C8:D8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Set the new values for C8 and D8
$ws.Range("C8").Value = 0.9925
$ws.Range("D8").Value = 0.9872

# Update the active selection to D8, matching the edited cell
$ws.Range("D8").Select() | Out-Null
